$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.194078947368421
$ws.Range("C2").Value = 0.555921052631579
$ws.Range("J2").Value = 0.01973684210526316
$ws.Range("P2").Value = 0.1447368421052632
$ws.Range("S2").Value = 0.08552631578947369
$ws.Range("B3").Value = 0.005376344086021506
$ws.Range("C3").Value = 0.02150537634408602
$ws.Range("J3").Value = 0.04838709677419355
$ws.Range("P3").Value = 0.7634408602150538
$ws.Range("S3").Value = 0.1612903225806452
$ws.Range("J4").Value = 0.09090909090909091
$ws.Range("O4").Value = 0.01818181818181818
$ws.Range("P4").Value = 0.5454545454545454
$ws.Range("S4").Value = 0.3454545454545455
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.6666666666666666
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("D6").Value = 0.01666666666666667
$ws.Range("F6").Value = 0.05
$ws.Range("J6").Value = 0.2125
$ws.Range("O6").Value = 0.025
$ws.Range("Q6").Value = 0.1958333333333333
$ws.Range("R6").Value = 0.07916666666666666
$ws.Range("S6").Value = 0.3541666666666667
$ws.Range("B7").Value = 0.1317365269461078
$ws.Range("D7").Value = 0.01197604790419162
$ws.Range("F7").Value = 0.04790419161676647
$ws.Range("J7").Value = 0.09580838323353294
$ws.Range("O7").Value = 0.01796407185628742
$ws.Range("Q7").Value = 0.2035928143712575
$ws.Range("R7").Value = 0.05389221556886228
$ws.Range("S7").Value = 0.437125748502994
$ws.Range("B8").Value = 0.0748792270531401
$ws.Range("D8").Value = 0.01449275362318841
$ws.Range("E8").Value = 0.002415458937198068
$ws.Range("F8").Value = 0.07246376811594203
$ws.Range("J8").Value = 0.108695652173913
$ws.Range("O8").Value = 0.01932367149758454
$ws.Range("Q8").Value = 0.2077294685990338
$ws.Range("R8").Value = 0.06763285024154589
$ws.Range("S8").Value = 0.4323671497584541
$ws.Range("B9").Value = 0.1105990783410138
$ws.Range("D9").Value = 0.02304147465437788
$ws.Range("F9").Value = 0.08755760368663594
$ws.Range("J9").Value = 0.08755760368663594
$ws.Range("O9").Value = 0.01382488479262673
$ws.Range("Q9").Value = 0.2258064516129032
$ws.Range("R9").Value = 0.08755760368663594
$ws.Range("S9").Value = 0.3640552995391705
$ws.Range("B10").Value = 0.1180223285486443
$ws.Range("D10").Value = 0.02791068580542265
$ws.Range("E10").Value = 0.001594896331738437
$ws.Range("F10").Value = 0.08373205741626795
$ws.Range("J10").Value = 0.1116427432216906
$ws.Range("O10").Value = 0.02312599681020734
$ws.Range("Q10").Value = 0.2129186602870813
$ws.Range("R10").Value = 0.07177033492822966
$ws.Range("S10").Value = 0.3492822966507177
$ws.Range("G11").Value = 0.1264822134387352
$ws.Range("J11").Value = 0.1067193675889328
$ws.Range("K11").Value = 0.1541501976284585
$ws.Range("L11").Value = 0.5889328063241107
$ws.Range("S11").Value = 0.02371541501976284
$ws.Range("G12").Value = 0.7672955974842768
$ws.Range("J12").Value = 0.1509433962264151
$ws.Range("K12").Value = 0.01257861635220126
$ws.Range("L12").Value = 0.05031446540880503
$ws.Range("S12").Value = 0.01886792452830189
$ws.Range("G13").Value = 0.5135135135135135
$ws.Range("J13").Value = 0.3783783783783784
$ws.Range("S13").Value = 0.1081081081081081
$ws.Range("F15").Value = 0.01587301587301587
$ws.Range("H15").Value = 0.1746031746031746
$ws.Range("I15").Value = 0.07539682539682539
$ws.Range("J15").Value = 0.3571428571428572
$ws.Range("K15").Value = 0.04761904761904762
$ws.Range("M15").Value = 0.0119047619047619
$ws.Range("O15").Value = 0.08333333333333333
$ws.Range("S15").Value = 0.2341269841269841
$ws.Range("F16").Value = 0.02538071065989848
$ws.Range("H16").Value = 0.1319796954314721
$ws.Range("I16").Value = 0.116751269035533
$ws.Range("J16").Value = 0.4213197969543147
$ws.Range("K16").Value = 0.116751269035533
$ws.Range("M16").Value = 0.02538071065989848
$ws.Range("O16").Value = 0.04568527918781726
$ws.Range("S16").Value = 0.116751269035533
$ws.Range("F17").Value = 0.01626016260162602
$ws.Range("H17").Value = 0.1788617886178862
$ws.Range("I17").Value = 0.1036585365853658
$ws.Range("J17").Value = 0.4573170731707317
$ws.Range("K17").Value = 0.07520325203252033
$ws.Range("M17").Value = 0.02032520325203252
$ws.Range("O17").Value = 0.06504065040650407
$ws.Range("S17").Value = 0.08333333333333333
$ws.Range("F18").Value = 0.01775147928994083
$ws.Range("H18").Value = 0.1893491124260355
$ws.Range("I18").Value = 0.1242603550295858
$ws.Range("J18").Value = 0.4378698224852071
$ws.Range("K18").Value = 0.07100591715976332
$ws.Range("M18").Value = 0.005917159763313609
$ws.Range("O18").Value = 0.106508875739645
$ws.Range("S18").Value = 0.04733727810650887
$ws.Range("F19").Value = 0.01663747810858144
$ws.Range("H19").Value = 0.1961471103327496
$ws.Range("I19").Value = 0.09194395796847636
$ws.Range("J19").Value = 0.3940455341506129
$ws.Range("K19").Value = 0.1129597197898424
$ws.Range("M19").Value = 0.01663747810858144
$ws.Range("N19").Value = 0.001751313485113835
$ws.Range("O19").Value = 0.08844133099824869
$ws.Range("S19").Value = 0.08143607705779335
